# CIV-6625 Update GA order template
# Remove the "Classification: Controlled" text-box shape (and its
# containing run) from the default footer.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)   # wdHeaderFooterPrimary - the default footer (footer2.xml)

if ($ftr.Exists -and $ftr.Shapes.Count -gt 0) {
    for ($i = $ftr.Shapes.Count; $i -ge 1; $i--) {
        $shape = $ftr.Shapes.Item($i)
        if ($shape.Name -eq "Text Box 4") {
            $shape.Delete()
        }
    }
}
